$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C, rows 2 through 42 hold the "Förändrad" (changed) date as a serial
# number. All of them move from 45713 (2025-02-25) to 45714 (2025-02-26),
# leaving formatting and all other cells untouched.
for ($row = 2; $row -le 42; $row++) {
    $ws.Cells.Item($row, 3).Value = 45714
}
